# edit.ps1
# Applies the changes described by the commit "Add files via upload":
#  1. Window view size (bookViews/workbookView windowWidth/windowHeight).
#  2. Updated "Source" footnote text (refreshed retrieval/publication dates).
#  3. Refreshed statistical values (oil/mineral/natural-resource rents and
#     related aggregates) in columns C:E for a number of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Workbook window size -------------------------------------------------
# Mirrors the author re-saving from a smaller Excel window
# (xWindow/yWindow=0, windowWidth 28800->19200, windowHeight 12490->10400).
try {
    $win = $excel.ActiveWindow
    $win.Width = 19200
    $win.Height = 10400
} catch {
    # Window chrome sizing may not be available in every host; ignore.
}

# --- 2. Updated source/footnote text -----------------------------------------
$ws.Range("A105").Value = 'Source : Organisation internationale du travail (récupéré le 14/11/2021), Division des statistiques des Nations Unies, Comptes nationaux (analyse des principaux agrégats, jeu de données téléchargé en décembre 2020), Indicateurs du développement dans le monde de la Banque mondiale (base de données et données publiées par les banques centrales, les agences nationales de statistique, et bureaux de pays de la Banque mondiale -22/04/2022).'

# --- 3. Refreshed data values -------------------------------------------------
$ws.Range("C13").Value = 4.24511870492843
$ws.Range("D13").Value = 1.0256180954427401
$ws.Range("E13").Value = 7.8760176205545198
$ws.Range("C23").Value = 9.2239408501254196
$ws.Range("D23").Value = 0.93080563470538003
$ws.Range("E23").Value = 15.041837465209699
$ws.Range("C38").Value = 0.64473634843504002
$ws.Range("D38").Value = 0.00052466462391999996
$ws.Range("E38").Value = 3.4898531808666799
$ws.Range("C45").Value = 7.3215173912970704
$ws.Range("D45").Value = 0.15644886761608001
$ws.Range("E45").Value = 8.5311452410706394
$ws.Range("C61").Value = 5.2818945746942001
$ws.Range("D61").Value = 0.070977671788809996
$ws.Range("E61").Value = 7.4865237466262897
$ws.Range("C62").Value = 5.2216836665081203
$ws.Range("D62").Value = 0.31636727112377
$ws.Range("E62").Value = 7.6403140161355996
$ws.Range("C63").Value = 1.53799932133146
$ws.Range("D63").Value = 0.21289374191068
$ws.Range("E63").Value = 2.3968920203694899
$ws.Range("C64").Value = 1.7736836525341499
$ws.Range("D64").Value = 0.44582365183315997
$ws.Range("E64").Value = 2.75337319723518
$ws.Range("C65").Value = 0.89492303171263998
$ws.Range("D65").Value = 0.30388181616974003
$ws.Range("E65").Value = 2.1546116573971399
$ws.Range("C66").Value = 1.7256311056939999
$ws.Range("D66").Value = 0.21816406658472001
$ws.Range("E66").Value = 2.6639606758106802
$ws.Range("C67").Value = 3.86105380973162
$ws.Range("D67").Value = 0.16784557187966001
$ws.Range("E67").Value = 6.0755247762998703
$ws.Range("C68").Value = 5.2781053276270002
$ws.Range("D68").Value = 0.12880948769287001
$ws.Range("E68").Value = 6.8848584335615302
$ws.Range("C69").Value = 0.0098146757952699993
$ws.Range("D69").Value = 0.00004643947943
$ws.Range("E69").Value = 2.57205177967192
$ws.Range("C70").Value = 14.8686072250135
$ws.Range("D70").Value = 0.52435270176505999
$ws.Range("E70").Value = 18.7344191769131
$ws.Range("C71").Value = 5.2818945746942001
$ws.Range("D71").Value = 0.070977671788809996
$ws.Range("E71").Value = 7.4865237466262897
$ws.Range("C72").Value = 0.85322440258807997
$ws.Range("D72").Value = 0.00070003148632999997
$ws.Range("E72").Value = 3.9202933658561601
$ws.Range("C73").Value = 3.416424340786
$ws.Range("D73").Value = 0.99593010995525
$ws.Range("E73").Value = 7.2048537359315299
$ws.Range("C74").Value = 11.125760654154799
$ws.Range("D74").Value = 0.3293511332772
$ws.Range("E74").Value = 12.623167783546901
$ws.Range("C75").Value = 20.646836147698099
$ws.Range("E75").Value = 23.889632314162199
$ws.Range("C76").Value = 0.90729207139557
$ws.Range("D76").Value = 0.17370901001059999
$ws.Range("E76").Value = 2.7886470056186701
$ws.Range("C77").Value = 1.9818701590583301
$ws.Range("D77").Value = 0.58838826316592996
$ws.Range("E77").Value = 3.2304196876906599
$ws.Range("C78").Value = 0.040160332233559998
$ws.Range("D78").Value = 0.03012799355049
$ws.Range("E78").Value = 0.19059421440464999
$ws.Range("C79").Value = 0.36314711207219003
$ws.Range("D79").Value = 0.11152221306819
$ws.Range("E79").Value = 0.67860633949164995
$ws.Range("C80").Value = 13.598975826458499
$ws.Range("D80").Value = 0.15081776776983
$ws.Range("E80").Value = 15.606151388553201
$ws.Range("C81").Value = 15.5253439794699
$ws.Range("D81").Value = 0.33605502046029001
$ws.Range("E81").Value = 18.309582978425802
$ws.Range("C82").Value = 1.55637400376145
$ws.Range("D82").Value = 0.38880000746870003
$ws.Range("E82").Value = 4.1550278239486902
$ws.Range("C83").Value = 0.45784838059269001
$ws.Range("D83").Value = 0.20338319870721
$ws.Range("E83").Value = 1.1681102423440199
$ws.Range("C84").Value = 1.2201164985295301
$ws.Range("D84").Value = 0.389526773682
$ws.Range("E84").Value = 6.1880612012389697
$ws.Range("C86").Value = 6.1996659926477697
$ws.Range("D86").Value = 0.11022997292483
$ws.Range("E86").Value = 7.7952395886384798
$ws.Range("C87").Value = 0.49907381428577002
$ws.Range("D87").Value = 0.43433652213343998
$ws.Range("E87").Value = 2.1858492446400701
$ws.Range("C88").Value = 4.9978301898136896
$ws.Range("D88").Value = 1.1120113589761
$ws.Range("E88").Value = 8.4143616566254895
$ws.Range("C89").Value = 2.2056835393462699
$ws.Range("D89").Value = 0.27170604001456
$ws.Range("E89").Value = 3.4272212523652801
$ws.Range("C90").Value = 1.37953764164502
$ws.Range("D90").Value = 0.10298128054202001
$ws.Range("E90").Value = 1.7244346023097801
$ws.Range("C91").Value = 4.1777892008159503
$ws.Range("D91").Value = 0.50230676143959996
$ws.Range("E91").Value = 8.2815665605199005
$ws.Range("C92").Value = 0.30778174802583003
$ws.Range("D92").Value = 0.14060517761255001
$ws.Range("E92").Value = 1.65738785692054
$ws.Range("E93").Value = 0.96646767838839998
$ws.Range("C94").Value = 0.47485830524744999
$ws.Range("D94").Value = 0.090686117353780002
$ws.Range("E94").Value = 1.3119703748365401
$ws.Range("C95").Value = 0.68471245342209996
$ws.Range("D95").Value = 0.19222733364474001
$ws.Range("E95").Value = 5.0090450209567798
$ws.Range("C96").Value = 6.8454985961850801
$ws.Range("D96").Value = 1.0424452023150601
$ws.Range("E96").Value = 10.461663234768
$ws.Range("C97").Value = 6.4939083086180904
$ws.Range("D97").Value = 0.26345531582738002
$ws.Range("E97").Value = 9.4327010113514103
$ws.Range("C98").Value = 5.6982155951968299
$ws.Range("D98").Value = 0.11983896057817001
$ws.Range("E98").Value = 6.7051665232767697
$ws.Range("C99").Value = 21.0021688369151
$ws.Range("D99").Value = 0.30032705968755002
$ws.Range("E99").Value = 22.9066097067589
